# Apply numeric corrections to the Durandal_Profits workbook (per-sheet Leve profit tables).
# Each assignment below sets a literal computed value for currentAveragePrice /
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ and LeveProfitNQ/HQ columns (H..N)
# on the row identified by its Leve Item ID, matching the scheduled runner's refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 2154.8125
$ws.Cells.Item(70, 9).Value = 1617
$ws.Cells.Item(70, 10).Value = 2477.5
$ws.Cells.Item(70, 11).Value = 4851
$ws.Cells.Item(70, 12).Value = 7432.5
$ws.Cells.Item(70, 13).Value = -4581
$ws.Cells.Item(70, 14).Value = -7972.5
$ws.Cells.Item(73, 8).Value = 2154.8125
$ws.Cells.Item(73, 9).Value = 1617
$ws.Cells.Item(73, 10).Value = 2477.5
$ws.Cells.Item(73, 11).Value = 4851
$ws.Cells.Item(73, 12).Value = 7432.5
$ws.Cells.Item(73, 13).Value = -3915
$ws.Cells.Item(73, 14).Value = -9304.5
$ws.Cells.Item(80, 8).Value = 2935.4333
$ws.Cells.Item(80, 9).Value = 845.53845
$ws.Cells.Item(80, 10).Value = 4533.5884
$ws.Cells.Item(80, 11).Value = 2536.61535
$ws.Cells.Item(80, 12).Value = 13600.7652
$ws.Cells.Item(80, 13).Value = -1538.61535
$ws.Cells.Item(80, 14).Value = -15596.7652
$ws.Cells.Item(83, 8).Value = 2935.4333
$ws.Cells.Item(83, 9).Value = 845.53845
$ws.Cells.Item(83, 10).Value = 4533.5884
$ws.Cells.Item(83, 11).Value = 7609.84605
$ws.Cells.Item(83, 12).Value = 40802.2956
$ws.Cells.Item(83, 13).Value = -2617.84605
$ws.Cells.Item(83, 14).Value = -50786.2956
$ws.Cells.Item(103, 8).Value = 83333930
$ws.Cells.Item(103, 9).Value = 674.5
$ws.Cells.Item(103, 10).Value = 125000550
$ws.Cells.Item(103, 11).Value = 2023.5
$ws.Cells.Item(103, 12).Value = 375001650
$ws.Cells.Item(103, 13).Value = -1437.5
$ws.Cells.Item(103, 14).Value = -375002822
$ws.Cells.Item(106, 8).Value = 2250.6
$ws.Cells.Item(106, 9).Value = 2197.0715
$ws.Cells.Item(106, 10).Value = 3000
$ws.Cells.Item(106, 11).Value = 2197.0715
$ws.Cells.Item(106, 12).Value = 3000
$ws.Cells.Item(106, 13).Value = -1566.0715
$ws.Cells.Item(106, 14).Value = -4262
$ws.Cells.Item(107, 8).Value = 758.55
$ws.Cells.Item(107, 9).Value = 821.82355
$ws.Cells.Item(107, 10).Value = 400
$ws.Cells.Item(107, 11).Value = 821.82355
$ws.Cells.Item(107, 12).Value = 400
$ws.Cells.Item(107, 13).Value = 1098.17645
$ws.Cells.Item(107, 14).Value = -4240
$ws.Cells.Item(113, 8).Value = 2972.1428
$ws.Cells.Item(113, 9).Value = 2752.5
$ws.Cells.Item(113, 10).Value = 3265
$ws.Cells.Item(113, 11).Value = 2752.5
$ws.Cells.Item(113, 12).Value = 3265
$ws.Cells.Item(113, 13).Value = 501.5
$ws.Cells.Item(113, 14).Value = -9773
$ws.Cells.Item(138, 8).Value = 3001.894
$ws.Cells.Item(138, 9).Value = 1629.35
$ws.Cells.Item(138, 10).Value = 3598.652
$ws.Cells.Item(138, 11).Value = 4888.049999999999
$ws.Cells.Item(138, 12).Value = 10795.956
$ws.Cells.Item(138, 13).Value = 251.9500000000007
$ws.Cells.Item(138, 14).Value = -21075.956

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 3815.8125
$ws.Cells.Item(2, 9).Value = 2970.1667
$ws.Cells.Item(2, 10).Value = 6352.75
$ws.Cells.Item(2, 11).Value = 2970.1667
$ws.Cells.Item(2, 12).Value = 6352.75
$ws.Cells.Item(2, 13).Value = -2857.1667
$ws.Cells.Item(2, 14).Value = -6578.75
$ws.Cells.Item(45, 8).Value = 2506.8333
$ws.Cells.Item(45, 9).Value = 2812.7334
$ws.Cells.Item(45, 10).Value = 2200.9333
$ws.Cells.Item(45, 11).Value = 2812.7334
$ws.Cells.Item(45, 12).Value = 2200.9333
$ws.Cells.Item(45, 13).Value = -2435.7334
$ws.Cells.Item(45, 14).Value = -2954.9333
$ws.Cells.Item(116, 8).Value = 3815.8125
$ws.Cells.Item(116, 9).Value = 2970.1667
$ws.Cells.Item(116, 10).Value = 6352.75
$ws.Cells.Item(116, 11).Value = 2970.1667
$ws.Cells.Item(116, 12).Value = 6352.75
$ws.Cells.Item(116, 13).Value = -676.1667000000002
$ws.Cells.Item(116, 14).Value = -10940.75
$ws.Cells.Item(122, 8).Value = 41250.117
$ws.Cells.Item(122, 9).Value = 2859
$ws.Cells.Item(122, 11).Value = 8577
$ws.Cells.Item(122, 13).Value = -6127
$ws.Cells.Item(132, 8).Value = 45502460
$ws.Cells.Item(132, 9).Value = 83335550
$ws.Cells.Item(132, 10).Value = 102742.8
$ws.Cells.Item(132, 11).Value = 250006650
$ws.Cells.Item(132, 12).Value = 308228.4
$ws.Cells.Item(132, 13).Value = -250004120
$ws.Cells.Item(132, 14).Value = -313288.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 3815.8125
$ws.Cells.Item(3, 9).Value = 2970.1667
$ws.Cells.Item(3, 10).Value = 6352.75
$ws.Cells.Item(3, 11).Value = 2970.1667
$ws.Cells.Item(3, 12).Value = 6352.75
$ws.Cells.Item(3, 13).Value = -2856.1667
$ws.Cells.Item(3, 14).Value = -6580.75
$ws.Cells.Item(80, 8).Value = 526.2381
$ws.Cells.Item(80, 10).Value = 343.57144
$ws.Cells.Item(80, 12).Value = 343.57144
$ws.Cells.Item(80, 14).Value = -2339.57144
$ws.Cells.Item(83, 8).Value = 526.2381
$ws.Cells.Item(83, 10).Value = 343.57144
$ws.Cells.Item(83, 12).Value = 1717.8572
$ws.Cells.Item(83, 14).Value = -11701.8572
$ws.Cells.Item(99, 8).Value = 1231.909
$ws.Cells.Item(99, 9).Value = 646.5294
$ws.Cells.Item(99, 10).Value = 3222.2
$ws.Cells.Item(99, 11).Value = 646.5294
$ws.Cells.Item(99, 12).Value = 3222.2
$ws.Cells.Item(99, 13).Value = 851.4706
$ws.Cells.Item(99, 14).Value = -6218.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1092.0952
$ws.Cells.Item(58, 9).Value = 1059.3334
$ws.Cells.Item(58, 10).Value = 1135.7778
$ws.Cells.Item(58, 11).Value = 1059.3334
$ws.Cells.Item(58, 12).Value = 1135.7778
$ws.Cells.Item(58, 13).Value = -856.3334
$ws.Cells.Item(58, 14).Value = -1541.7778
$ws.Cells.Item(94, 8).Value = 4009.158
$ws.Cells.Item(94, 10).Value = 5061.5713
$ws.Cells.Item(94, 12).Value = 5061.5713
$ws.Cells.Item(94, 14).Value = -5963.5713
$ws.Cells.Item(132, 8).Value = 46182.61
$ws.Cells.Item(132, 9).Value = 1887.5
$ws.Cells.Item(132, 10).Value = 147428.58
$ws.Cells.Item(132, 11).Value = 5662.5
$ws.Cells.Item(132, 12).Value = 442285.74
$ws.Cells.Item(132, 13).Value = -3132.5
$ws.Cells.Item(132, 14).Value = -447345.74
$ws.Cells.Item(136, 8).Value = 1092.0952
$ws.Cells.Item(136, 9).Value = 1059.3334
$ws.Cells.Item(136, 10).Value = 1135.7778
$ws.Cells.Item(136, 11).Value = 3178.0002
$ws.Cells.Item(136, 12).Value = 3407.3334
$ws.Cells.Item(136, 13).Value = -628.0001999999999
$ws.Cells.Item(136, 14).Value = -8507.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 1215.3846
$ws.Cells.Item(68, 9).Value = 985.7143
$ws.Cells.Item(68, 10).Value = 1483.3334
$ws.Cells.Item(68, 11).Value = 2957.1429
$ws.Cells.Item(68, 12).Value = 4450.0002
$ws.Cells.Item(68, 13).Value = -2146.1429
$ws.Cells.Item(68, 14).Value = -6072.0002
$ws.Cells.Item(71, 8).Value = 1215.3846
$ws.Cells.Item(71, 9).Value = 985.7143
$ws.Cells.Item(71, 10).Value = 1483.3334
$ws.Cells.Item(71, 11).Value = 8871.4287
$ws.Cells.Item(71, 12).Value = 13350.0006
$ws.Cells.Item(71, 13).Value = -4815.4287
$ws.Cells.Item(71, 14).Value = -21462.0006
$ws.Cells.Item(107, 8).Value = 256.06818
$ws.Cells.Item(107, 9).Value = 242
$ws.Cells.Item(107, 10).Value = 259.6857
$ws.Cells.Item(107, 11).Value = 726
$ws.Cells.Item(107, 12).Value = 779.0571
$ws.Cells.Item(107, 13).Value = 1194
$ws.Cells.Item(107, 14).Value = -4619.0571
$ws.Cells.Item(113, 8).Value = 1020.2639
$ws.Cells.Item(113, 9).Value = 1003
$ws.Cells.Item(113, 10).Value = 1020.507
$ws.Cells.Item(113, 11).Value = 3009
$ws.Cells.Item(113, 12).Value = 3061.521
$ws.Cells.Item(113, 13).Value = -839
$ws.Cells.Item(113, 14).Value = -7401.521

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 1524.85
$ws.Cells.Item(122, 9).Value = 1526.4546
$ws.Cells.Item(122, 10).Value = 1522.8889
$ws.Cells.Item(122, 11).Value = 4579.3638
$ws.Cells.Item(122, 12).Value = 4568.6667
$ws.Cells.Item(122, 13).Value = -2129.3638
$ws.Cells.Item(122, 14).Value = -9468.6667
$ws.Cells.Item(126, 8).Value = 10424328
$ws.Cells.Item(126, 9).Value = 13686.5
$ws.Cells.Item(126, 10).Value = 20834970
$ws.Cells.Item(126, 11).Value = 41059.5
$ws.Cells.Item(126, 12).Value = 62504910
$ws.Cells.Item(126, 13).Value = -38589.5
$ws.Cells.Item(126, 14).Value = -62509850
$ws.Cells.Item(132, 8).Value = 306916.06
$ws.Cells.Item(132, 9).Value = 48673.145
$ws.Cells.Item(132, 10).Value = 668456.1
$ws.Cells.Item(132, 11).Value = 146019.435
$ws.Cells.Item(132, 12).Value = 2005368.3
$ws.Cells.Item(132, 13).Value = -143489.435
$ws.Cells.Item(132, 14).Value = -2010428.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 32482.312
$ws.Cells.Item(132, 9).Value = 56812
$ws.Cells.Item(132, 10).Value = 1201.2858
$ws.Cells.Item(132, 11).Value = 170436
$ws.Cells.Item(132, 12).Value = 3603.8574
$ws.Cells.Item(132, 13).Value = -167906
$ws.Cells.Item(132, 14).Value = -8663.857400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 58825220
$ws.Cells.Item(126, 9).Value = 125001464
$ws.Cells.Item(126, 10).Value = 1893.8889
$ws.Cells.Item(126, 11).Value = 375004392
$ws.Cells.Item(126, 12).Value = 5681.6667
$ws.Cells.Item(126, 13).Value = -375001922
$ws.Cells.Item(126, 14).Value = -10621.6667
$ws.Cells.Item(132, 8).Value = 152825730
$ws.Cells.Item(132, 9).Value = 226001890
$ws.Cells.Item(132, 10).Value = 6473414
$ws.Cells.Item(132, 11).Value = 678005670
$ws.Cells.Item(132, 12).Value = 19420242
$ws.Cells.Item(132, 13).Value = -678003140
$ws.Cells.Item(132, 14).Value = -19425302
$ws.Cells.Item(136, 8).Value = 36014.586
$ws.Cells.Item(136, 9).Value = 56391
$ws.Cells.Item(136, 10).Value = 2671.3635
$ws.Cells.Item(136, 11).Value = 169173
$ws.Cells.Item(136, 12).Value = 8014.0905
$ws.Cells.Item(136, 13).Value = -166623
$ws.Cells.Item(136, 14).Value = -13114.0905
